$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (Coin name) updates ---
$bUpdates = @{
    48 = 'Filecoin'
    49 = 'BabyDogeCoin'
}
foreach ($row in $bUpdates.Keys) {
    $ws.Cells.Item($row, 2).Value = $bUpdates[$row]
}

# --- Column C (Link) updates ---
$cUpdates = @{
    48 = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    49 = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
}
foreach ($row in $cUpdates.Keys) {
    $ws.Cells.Item($row, 3).Value = $cUpdates[$row]
}

# --- Column D (Price) updates; force text so numeric-looking strings are not coerced ---
$dUpdates = @{
    2 = '65.753.28'
    3 = '2.499.17'
    5 = '577.06'
    6 = '166.32'
    8 = '0.518'
    9 = '2.498.00'
    10 = '0.136'
    13 = '5.09'
    14 = '26.17'
    15 = '2.957.98'
    16 = '0.0000174'
    17 = '65.632.77'
    18 = '2.461.33'
    19 = '11.18'
    20 = '7.56'
    21 = '342.48'
    22 = '4.18'
    23 = '4.52'
    25 = '1.92'
    26 = '68.84'
    27 = '9.89'
    29 = '2.631.39'
    30 = '0.0₃0971'
    31 = '8.10'
    32 = '518.98'
    33 = '1.30'
    34 = '1.80'
    37 = '157.09'
    39 = '18.41'
    41 = '0.353'
    43 = '5.02'
    45 = '2.42'
    46 = '146.06'
    47 = '0.552'
    48 = '3.67'
    49 = '0.0₆0273'
    51 = '0.0749'
}
foreach ($row in $dUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $dUpdates[$row]
    $cell.Style = "Normal"
}

# --- Column E (Volume 1h) updates ---
$eUpdates = @{
    2 = '  -2.36%  '
    3 = '  -4.78%  '
    4 = '  +0.02%  '
    5 = '  -3.07%  '
    6 = '  -0.33%  '
    7 = '  +0.10%  '
    8 = '  -2.85%  '
    9 = '  -4.82%  '
    10 = '  -1.77%  '
    11 = '  -0.29%  '
    12 = '  -4.45%  '
    13 = '  -2.85%  '
    14 = '  -5.36%  '
    15 = '  -4.64%  '
    16 = '  -4.60%  '
    17 = '  -2.40%  '
    18 = '  -6.07%  '
    19 = '  -7.29%  '
    20 = '  -5.10%  '
    21 = '  -4.31%  '
    22 = '  -3.37%  '
    23 = '  -3.08%  '
    24 = '  +0.01%  '
    25 = '  -0.68%  '
    26 = '  -1.33%  '
    27 = '  -3.60%  '
    28 = '  +0.29%  '
    29 = '  -4.65%  '
    30 = '  -3.27%  '
    31 = '  +2.31%  '
    32 = '  -5.01%  '
    33 = '  -3.49%  '
    34 = '  -5.31%  '
    35 = '  -4.54%  '
    36 = '  +0.00%  '
    37 = '  +0.02%  '
    38 = '  -4.42%  '
    39 = '  -3.14%  '
    40 = '  +0.58%  '
    41 = '  -3.74%  '
    42 = '  -3.54%  '
    43 = '  -3.80%  '
    44 = '  -0.05%  '
    45 = '  -0.09%  '
    46 = '  -3.97%  '
    47 = '  -4.92%  '
    48 = '  -2.75%  '
    49 = '  -8.16%  '
    50 = '  +0.20%  '
    51 = '  -2.62%  '
}
foreach ($row in $eUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $eUpdates[$row]
}
